$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("K3").Value = 2.5
$ws.Range("L3").Value = 8.5
$ws.Range("O3").Value = 1.22
$ws.Range("P3").Value = 4.33
$ws.Range("Q3").Value = 1.67
$ws.Range("R3").Value = 2.2
$ws.Range("S3").Value = 1.3
$ws.Range("T3").Value = 3.4
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62
$ws.Range("W3").Value = 6.5
$ws.Range("Y3").Value = 8.5
$ws.Range("Z3").Value = 8
$ws.Range("AA3").Value = 12
$ws.Range("AF3").Value = 101
$ws.Range("AH3").Value = 17
$ws.Range("AJ3").Value = 21
$ws.Range("AL3").Value = 51
$ws.Range("AP3").Value = 19
$ws.Range("AQ3").Value = 17
$ws.Range("AT3").Value = 3.4
$ws.Range("G5").Value = 1.62
$ws.Range("M5").Value = 1.03
$ws.Range("N5").Value = 15
$ws.Range("Q5").Value = 1.67
$ws.Range("R5").Value = 2.15
$ws.Range("Z5").Value = 12
$ws.Range("AO5").Value = 8
$ws.Range("AQ5").Value = 23
$ws.Range("AW5").Value = 7
$ws.Range("AX5").Value = 26
$ws.Range("BB5").Value = 201
$ws.Range("I6").Value = 2.63
$ws.Range("Q6").Value = 1.73
$ws.Range("R6").Value = 2.08
$ws.Range("U6").Value = 1.57
$ws.Range("V6").Value = 2.25
$ws.Range("W6").Value = 11
$ws.Range("AB6").Value = 23
$ws.Range("AL6").Value = 19
$ws.Range("S7").Value = 1.23
$ws.Range("T7").Value = 3.9
$ws.Range("U7").Value = 2
$ws.Range("V7").Value = 1.65
$ws.Range("G10").Value = 1.75
$ws.Range("K10").Value = 2.5
$ws.Range("L10").Value = 4.33
$ws.Range("O10").Value = 1.14
$ws.Range("P10").Value = 5.5
$ws.Range("S10").Value = 1.25
$ws.Range("T10").Value = 3.75
$ws.Range("W10").Value = 11
$ws.Range("X10").Value = 10
$ws.Range("AH10").Value = 19
$ws.Range("AP10").Value = 15
$ws.Range("AT10").Value = 3.75
$ws.Range("AW10").Value = 6.5
$ws.Range("AZ10").Value = 67
$ws.Range("G11").Value = 1.73
$ws.Range("H11").Value = 4.1
$ws.Range("I11").Value = 4.1
$ws.Range("J11").Value = 2.3
$ws.Range("K11").Value = 2.3
$ws.Range("L11").Value = 4.5
$ws.Range("U11").Value = 1.67
$ws.Range("V11").Value = 2.1
$ws.Range("W11").Value = 8.5
$ws.Range("X11").Value = 9
$ws.Range("Z11").Value = 15
$ws.Range("AD11").Value = 7.5
$ws.Range("AF11").Value = 41
$ws.Range("AJ11").Value = 13
$ws.Range("AK11").Value = 41
$ws.Range("AN11").Value = 4
$ws.Range("AO11").Value = 9
$ws.Range("AQ11").Value = 29
$ws.Range("AU11").Value = 7.5
$ws.Range("AY11").Value = 26
$ws.Range("AZ11").Value = 67
$ws.Range("BA11").Value = 81
$ws.Range("H12").Value = 5.5
$ws.Range("M12").Value = 1.05
$ws.Range("N12").Value = 11
$ws.Range("J21").Value = 1.85
$ws.Range("L21").Value = 6
$ws.Range("P21").Value = 4.6
$ws.Range("T21").Value = 3.35
$ws.Range("W21").Value = 9.25
$ws.Range("X21").Value = 8.25
$ws.Range("Z21").Value = 10.5
$ws.Range("AA21").Value = 10.25
$ws.Range("AB21").Value = 19.5
$ws.Range("AE21").Value = 16
$ws.Range("AG21").Value = 350
$ws.Range("AH21").Value = 21
$ws.Range("AI21").Value = 45
$ws.Range("AL21").Value = 65
$ws.Range("AM21").Value = 50
$ws.Range("AO21").Value = 6.4
$ws.Range("AP21").Value = 13.5
$ws.Range("AQ21").Value = 17
$ws.Range("AR21").Value = 35
$ws.Range("AT21").Value = 3.35
$ws.Range("AU21").Value = 7.4
$ws.Range("AV21").Value = 55
$ws.Range("AX21").Value = 37
$ws.Range("AY21").Value = 35
$ws.Range("AZ21").Value = 250
$ws.Range("BA21").Value = 200
$ws.Range("BB21").Value = 350
$ws.Range("I23").Value = 2.75
$ws.Range("O23").Value = 1.35
$ws.Range("P23").Value = 2.7
$ws.Range("Q23").Value = 2.02
$ws.Range("AT23").Value = 2.47
$ws.Range("G38").Value = 3.7
$ws.Range("H38").Value = 3.25
$ws.Range("I38").Value = 1.91
$ws.Range("J38").Value = 4.3
$ws.Range("L38").Value = 2.57
$ws.Range("Q38").Value = 2.2
$ws.Range("R38").Value = 1.62
$ws.Range("U38").Value = 1.98
$ws.Range("W38").Value = 9.25
$ws.Range("X38").Value = 18.5
$ws.Range("AA38").Value = 40
$ws.Range("AD38").Value = 6.5
$ws.Range("AE38").Value = 17.5
$ws.Range("AG38").Value = 1000
$ws.Range("AH38").Value = 6.1
$ws.Range("AI38").Value = 8.25
$ws.Range("AK38").Value = 16
$ws.Range("AL38").Value = 17
$ws.Range("AM38").Value = 35
$ws.Range("AT38").Value = 2.52
$ws.Range("AU38").Value = 7.9
$ws.Range("AV38").Value = 90
$ws.Range("AW38").Value = 3.7
$ws.Range("AX38").Value = 10
$ws.Range("AY38").Value = 22
$ws.Range("BA38").Value = 90
